# The `reviews_count` column (originally column E) was removed from the
# report. Deleting the whole column shifts every column to its right
# (reviews_average, latitude, longitude, is_permanently_closed, gmaps_link,
# latest_review_date) one position to the left, which matches the target
# diff (old F:K -> new E:J) and shrinks the used range from A1:K21 to A1:J21.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns("E").Delete()
